# kadın erkek bar plot & pie chart
# Rename the age-group category labels (column A, rows 2-97) on both sheets
# from the old "Yas_16_29" style codes to human-readable Turkish labels, and
# tidy up a couple of leftover view/format bits (selection + column A width
# + the header cell A1's alignment on the first sheet).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Kadın Yaş Dağılımı"
$ws2 = $wb.Worksheets.Item(2)   # "Erkek Yaş Dağılımı"

foreach ($ws in @($ws1, $ws2)) {
    $ws.Range("A2:A25").Value  = "Yaş Aralığı: 16-29"
    $ws.Range("A26:A49").Value = "Yaş Aralığı: 30-44"
    $ws.Range("A50:A73").Value = "Yaş Aralığı: 45-54"
    $ws.Range("A74:A97").Value = "Yaş Aralığı: 55 ve 60 üstü"
}

# --- Sheet 1 ("Kadın Yaş Dağılımı") cosmetic touch-ups -------------------

# A1's alignment used to be "vertical center" (style index 2); it should
# match the rest of the header row, i.e. "horizontal center" (style index
# 1, same as B1). Copy B1's formatting onto A1 without touching its value.
$ws1.Range("B1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws1.Range("A1").Value = "Yas_Grubu"

# Column A is a bit wider now.
$ws1.Columns.Item(1).ColumnWidth = 26.5

# The remembered selection on sheet 1 moved from E20 to B2.
[void]$ws1.Range("B2").Select()

# --- Sheet 2 ("Erkek Yaş Dağılımı") cosmetic touch-ups --------------------

# Re-activate sheet 2 (it is the tab shown when the file is opened) and
# reset its remembered selection back to the top-left cell.
[void]$ws2.Activate()
[void]$ws2.Range("A1").Select()
$excel.CutCopyMode = 0
